$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Write stories" row (row 11): was assigned to "Nate, Xavier, Ali" with no
# time estimate; now it's broken out with its own time/ risk/ completion
# data, matching the other finished tasks (Set up group GitHub account /
# Set up Slack channel).
$ws.Range("B11").Value = "15 min"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Nate"
$ws.Range("E11").Value = "20 min"
$ws.Range("F11").Value = 1
$ws.Range("F11").NumberFormat = "0%"

# Move the active selection to reflect where editing left off.
$ws.Range("F15").Select()
